$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new column before the current column A ("Name"),
#    shifting Name/Description/Price/Category one column to the right.
$ws.Columns.Item(1).Insert()

# --- Header row -----------------------------------------------------
$ws.Range("A1").Value = "Code"
$ws.Range("E1").Value = "Category"

# Give the "Price" header (now column D) a numeric-style number format,
# keeping the existing bold/border/left-aligned header look.
$ws.Range("D1").NumberFormat = "#,##0"

# --- Price column (now D) values ------------------------------------
# Add a border + right alignment + explicit font + number format to the
# price cells.
$ws.Range("D2:D4").Font.Name = "Calibri"
$ws.Range("D2:D4").Borders.Item(1).LineStyle = 1
$ws.Range("D2:D4").HorizontalAlignment = -4152
$ws.Range("D2:D4").NumberFormat = "#,##0"

# Column-level default formatting for the Price column.
$ws.Columns.Item(4).NumberFormat = "#,##0"
$ws.Columns.Item(4).HorizontalAlignment = -4152

# --- Code column (new A) ---------------------------------------------
# Bold header, text format, no border.
$ws.Range("A1").NumberFormat = "@"
$ws.Range("A1").Font.Bold = $true

# Text-formatted code values (protect leading zeroes).
$ws.Range("A2:A4").NumberFormat = "@"
$ws.Range("A2").Value = "0001"
$ws.Range("A3").Value = "0002"
$ws.Range("A4").Value = "0003"

$ws.Columns.Item(1).NumberFormat = "@"
$ws.Columns.Item(1).ColumnWidth = 11.71875

# --- Category column (new E) values -----------------------------------
$ws.Range("E2").Value = "Hambúrguer"
$ws.Range("E3").Value = "Hambúrguer"
$ws.Range("E4").Value = "Hambúrguer"

# --- Row heights for the data rows -----------------------------------
$ws.Rows.Item(2).RowHeight = 17.25
$ws.Rows.Item(3).RowHeight = 17.25
$ws.Rows.Item(4).RowHeight = 17.25

# --- Final selection ---------------------------------------------------
$ws.Range("A5").Select()
